$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at H; existing H:K shift to I:L, inheriting column G's style.
$ws.Columns.Item(8).Insert()

# Split "Ministry Course Code and Level" header into two headers.
$ws.Range("G1").Value = "Ministry Course Code"

# Split the "ENST 12" values in column G into code (G) / level (H, numeric).
$ws.Range("G2").Value = "ENST"
$ws.Range("G3").Value = "ENST"
$ws.Range("G4").Value = "ENST"

$ws.Range("H1").Value = "Ministry Course Level"
$ws.Range("H2").Value = 12
$ws.Range("H3").Value = 12
$ws.Range("H4").Value = 12

# Mirror the selection left on the new columns after the split.
$ws.Range("G1:H1048576").Select() | Out-Null
